# "Điều chỉnh tên bài và tên tệp" - adjust the lesson title/subtitle textbox
# on slide 1: reposition/resize it, rename "Chương 5" -> "Bài 5", and add the
# missing space after "IIR.C11." on the second line.

$p   = $ppt.ActivePresentation
$s   = $p.Slides.Item(1)
$shp = $s.Shapes("Rectangle 3")

# --- Move / resize the subtitle placeholder -------------------------------
# Shape.Left/Top/Width/Height are in points; the XML stores EMU (1 pt = 12700 EMU).
$shp.Left   = 611560  / 12700
$shp.Top    = 3429000 / 12700
$shp.Width  = 7704856 / 12700
$shp.Height = 1968624 / 12700

$tr = $shp.TextFrame.TextRange

# --- Line 1: "Chương 5.Mô " -> "Bài 5. Mô " -------------------------------
# Clear the leading word "Chương" first, then retype the remainder so the
# replacement adopts the following run's formatting instead of the deleted
# run's formatting.
$tr.Characters(1, 6).Text = ""
$shp.TextFrame.TextRange.Characters(1, 6).Text = "Bài 5. Mô "

# --- Line 2: "IIR.C11.Probabilistic " -> "IIR.C11" + ". Probabilistic " ---
# Remove the "." joining "IIR.C11" and "Probabilistic", then reinsert ". "
# in front of "Probabilistic " so the text splits into two runs.
$tr2 = $shp.TextFrame.TextRange
$tr2.Characters(40, 1).Text = ""
$shp.TextFrame.TextRange.Characters(40, 14).InsertBefore(". ")
